$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 40 new rows before row 4414, shifting existing rows (old 4414-4421) down to 4454-4461
$ws.Range("A4414:A4453").EntireRow.Insert()

$ws.Cells.Item(4414, 1).Value = "Basica"
$ws.Cells.Item(4414, 2).Value = "'2023-08-09"
$ws.Cells.Item(4414, 3).Value = 27
$ws.Cells.Item(4415, 1).Value = "Media"
$ws.Cells.Item(4415, 2).Value = "'2023-08-09"
$ws.Cells.Item(4415, 3).Value = 32
$ws.Cells.Item(4416, 1).Value = "UTI"
$ws.Cells.Item(4416, 2).Value = "'2023-08-09"
$ws.Cells.Item(4416, 3).Value = 7
$ws.Cells.Item(4417, 1).Value = "UCI"
$ws.Cells.Item(4417, 2).Value = "'2023-08-09"
$ws.Cells.Item(4417, 3).Value = 9
$ws.Cells.Item(4418, 1).Value = "Basica"
$ws.Cells.Item(4418, 2).Value = "'2023-08-10"
$ws.Cells.Item(4418, 3).Value = 23
$ws.Cells.Item(4419, 1).Value = "Media"
$ws.Cells.Item(4419, 2).Value = "'2023-08-10"
$ws.Cells.Item(4419, 3).Value = 37
$ws.Cells.Item(4420, 1).Value = "UTI"
$ws.Cells.Item(4420, 2).Value = "'2023-08-10"
$ws.Cells.Item(4420, 3).Value = 5
$ws.Cells.Item(4421, 1).Value = "UCI"
$ws.Cells.Item(4421, 2).Value = "'2023-08-10"
$ws.Cells.Item(4421, 3).Value = 11
$ws.Cells.Item(4422, 1).Value = "Basica"
$ws.Cells.Item(4422, 2).Value = "'2023-08-11"
$ws.Cells.Item(4422, 3).Value = 25
$ws.Cells.Item(4423, 1).Value = "Media"
$ws.Cells.Item(4423, 2).Value = "'2023-08-11"
$ws.Cells.Item(4423, 3).Value = 37
$ws.Cells.Item(4424, 1).Value = "UTI"
$ws.Cells.Item(4424, 2).Value = "'2023-08-11"
$ws.Cells.Item(4424, 3).Value = 7
$ws.Cells.Item(4425, 1).Value = "UCI"
$ws.Cells.Item(4425, 2).Value = "'2023-08-11"
$ws.Cells.Item(4425, 3).Value = 10
$ws.Cells.Item(4426, 1).Value = "Basica"
$ws.Cells.Item(4426, 2).Value = "'2023-08-12"
$ws.Cells.Item(4426, 3).Value = 22
$ws.Cells.Item(4427, 1).Value = "Media"
$ws.Cells.Item(4427, 2).Value = "'2023-08-12"
$ws.Cells.Item(4427, 3).Value = 26
$ws.Cells.Item(4428, 1).Value = "UTI"
$ws.Cells.Item(4428, 2).Value = "'2023-08-12"
$ws.Cells.Item(4428, 3).Value = 8
$ws.Cells.Item(4429, 1).Value = "UCI"
$ws.Cells.Item(4429, 2).Value = "'2023-08-12"
$ws.Cells.Item(4429, 3).Value = 12
$ws.Cells.Item(4430, 1).Value = "Basica"
$ws.Cells.Item(4430, 2).Value = "'2023-08-14"
$ws.Cells.Item(4430, 3).Value = 22
$ws.Cells.Item(4431, 1).Value = "Media"
$ws.Cells.Item(4431, 2).Value = "'2023-08-14"
$ws.Cells.Item(4431, 3).Value = 27
$ws.Cells.Item(4432, 1).Value = "UTI"
$ws.Cells.Item(4432, 2).Value = "'2023-08-14"
$ws.Cells.Item(4432, 3).Value = 7
$ws.Cells.Item(4433, 1).Value = "UCI"
$ws.Cells.Item(4433, 2).Value = "'2023-08-14"
$ws.Cells.Item(4433, 3).Value = 11
$ws.Cells.Item(4434, 1).Value = "Basica"
$ws.Cells.Item(4434, 2).Value = "'2023-08-15"
$ws.Cells.Item(4434, 3).Value = 24
$ws.Cells.Item(4435, 1).Value = "Media"
$ws.Cells.Item(4435, 2).Value = "'2023-08-15"
$ws.Cells.Item(4435, 3).Value = 30
$ws.Cells.Item(4436, 1).Value = "UTI"
$ws.Cells.Item(4436, 2).Value = "'2023-08-15"
$ws.Cells.Item(4436, 3).Value = 11
$ws.Cells.Item(4437, 1).Value = "UCI"
$ws.Cells.Item(4437, 2).Value = "'2023-08-15"
$ws.Cells.Item(4437, 3).Value = 9
$ws.Cells.Item(4438, 1).Value = "Basica"
$ws.Cells.Item(4438, 2).Value = "'2023-08-16"
$ws.Cells.Item(4438, 3).Value = 19
$ws.Cells.Item(4439, 1).Value = "Media"
$ws.Cells.Item(4439, 2).Value = "'2023-08-16"
$ws.Cells.Item(4439, 3).Value = 22
$ws.Cells.Item(4440, 1).Value = "UTI"
$ws.Cells.Item(4440, 2).Value = "'2023-08-16"
$ws.Cells.Item(4440, 3).Value = 11
$ws.Cells.Item(4441, 1).Value = "UCI"
$ws.Cells.Item(4441, 2).Value = "'2023-08-16"
$ws.Cells.Item(4441, 3).Value = 12
$ws.Cells.Item(4442, 1).Value = "Basica"
$ws.Cells.Item(4442, 2).Value = "'2023-08-17"
$ws.Cells.Item(4442, 3).Value = 29
$ws.Cells.Item(4443, 1).Value = "Media"
$ws.Cells.Item(4443, 2).Value = "'2023-08-17"
$ws.Cells.Item(4443, 3).Value = 19
$ws.Cells.Item(4444, 1).Value = "UTI"
$ws.Cells.Item(4444, 2).Value = "'2023-08-17"
$ws.Cells.Item(4444, 3).Value = 13
$ws.Cells.Item(4445, 1).Value = "UCI"
$ws.Cells.Item(4445, 2).Value = "'2023-08-17"
$ws.Cells.Item(4445, 3).Value = 12
$ws.Cells.Item(4446, 1).Value = "Basica"
$ws.Cells.Item(4446, 2).Value = "'2023-08-18"
$ws.Cells.Item(4446, 3).Value = 21
$ws.Cells.Item(4447, 1).Value = "Media"
$ws.Cells.Item(4447, 2).Value = "'2023-08-18"
$ws.Cells.Item(4447, 3).Value = 20
$ws.Cells.Item(4448, 1).Value = "UTI"
$ws.Cells.Item(4448, 2).Value = "'2023-08-18"
$ws.Cells.Item(4448, 3).Value = 11
$ws.Cells.Item(4449, 1).Value = "UCI"
$ws.Cells.Item(4449, 2).Value = "'2023-08-18"
$ws.Cells.Item(4449, 3).Value = 10
$ws.Cells.Item(4450, 1).Value = "Basica"
$ws.Cells.Item(4450, 2).Value = "'2023-08-19"
$ws.Cells.Item(4450, 3).Value = 35
$ws.Cells.Item(4451, 1).Value = "Media"
$ws.Cells.Item(4451, 2).Value = "'2023-08-19"
$ws.Cells.Item(4451, 3).Value = 24
$ws.Cells.Item(4452, 1).Value = "UTI"
$ws.Cells.Item(4452, 2).Value = "'2023-08-19"
$ws.Cells.Item(4452, 3).Value = 9
$ws.Cells.Item(4453, 1).Value = "UCI"
$ws.Cells.Item(4453, 2).Value = "'2023-08-19"
$ws.Cells.Item(4453, 3).Value = 10

# Re-apply the plain/default formatting (matches the style of surrounding data rows)
# so the text-looking dates above don't retain a quote-prefix style.
$ws.Range("A4413:C4413").Copy()
$ws.Range("A4414:C4453").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "New dimension:" $ws.UsedRange.Address()
Write-Host "B4414:" $ws.Cells.Item(4414,2).Text
Write-Host "C4453:" $ws.Cells.Item(4453,3).Text
Write-Host "B4454 (shifted old data):" $ws.Cells.Item(4454,2).Text
Write-Host "C4461 (shifted old data):" $ws.Cells.Item(4461,3).Text
